$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update 2D training schedule values (row 2)
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 5
$ws.Cells.Item(2, 5).Value = 4
$ws.Cells.Item(2, 6).Value = 3
$ws.Cells.Item(2, 7).Value = 3
$ws.Cells.Item(2, 8).Value = 34

# row 3
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 6
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 5
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 12

# row 4
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 2
$ws.Cells.Item(4, 5).Value = 8

# row 5
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = 7
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(5, 8).Value = 23

# row 6
$ws.Cells.Item(6, 2).Value = 3
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 5
$ws.Cells.Item(6, 5).Value = 4

# Update the active selected cell from I2 to I1
$ws.Range("I1").Select()

$wb.Save()
